# "Script TMPro update functioning"
#
# The XML-mapped "Progressions" table on Taul1 (A1:K21) leaves the
# Advancement1..Advancement10 cells blank whenever a level has fewer than
# ten advancements. The update script now back-fills every still-empty
# cell inside the table body (rows 2-21, columns B-K) with the literal
# text "none" so downstream consumers never see a blank cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

$firstRow = 2
$lastRow  = 21
$firstCol = 2   # B
$lastCol  = 11  # K

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $null) {
            $cell.Value = "none"
        }
    }
}
